# Rename the second worksheet (strategy_id-6002 -> strategy_id-6005)
$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "strategy_id-6005"

# Clear out column C data (keep header in row 1) on both sheets
foreach ($ws in $wb.Worksheets) {
    $usedRange = $ws.UsedRange
    $lastRow = $usedRange.Rows.Count + $usedRange.Row - 1
    if ($lastRow -ge 2) {
        $ws.Range("C2:C" + $lastRow).ClearContents()
    }
}
